$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nov")

# Add new header and values in column G
$ws.Range("G1").Value = "newMessage"
$ws.Range("G2").Value = "new"
$ws.Range("G3").Value = "new"
$ws.Range("G4").Value = "new"

# Update selection to match the diff (G5)
$ws.Range("G5").Select()
